$d = $word.ActiveDocument

# XML fragments for each of the 50 new paragraphs to insert, in order,
# right before the final (already existing) empty trailing paragraph.
$paraXml = @(
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>------------------------------------------------------------------------------------------------------------------------</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">S9: install ansible, </w:t></w:r><w:r><w:rPr/><w:t>create scripts and run them – check for error in script</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>sudo apt install ansible -y</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Set the ip addresses, ansible_user and </w:t></w:r><w:r><w:rPr/><w:t>private key for the user in the hosts file</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>vi hosts</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>[flask_vm]</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>&lt;ip of machine&gt; ansible_user=azureuser ansible_ssh_private_key_file=~/.ssh/azure_id_rsa</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>vi deploy_flask.yml</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>- name: Deploy Flask App</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr/><w:t>hosts: flask_app_vm</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr/><w:t>become: yes</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr/><w:t>tasks:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>- name: Install Python3 and Pip</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr/><w:t>apt:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>update_cache: yes</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>- name: Install Python3 and Pip</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr/><w:t>apt:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>name:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">          </w:t></w:r><w:r><w:rPr/><w:t>- python3</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">          </w:t></w:r><w:r><w:rPr/><w:t>- python3-pip</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>state: present</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>- name: Install Flask</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr/><w:t>pip:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>name: flask</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>- name: Copy Flask App to VM</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr/><w:t>copy:</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>src: ../flask_app/app.py</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr/><w:t>dest: /home/azureuser/flask_app/</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr/><w:t>- name: Start Flask App</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr/><w:t>shell: python3 /home/azureuser/flask_app/app.py</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>ansible-playbook -i hosts deploy_flask.yml</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>Test on &lt;public ip:5000&gt;</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>',
)

$lastIdx = $d.Paragraphs.Count
$anchor = $d.Paragraphs($lastIdx).Range

# Create N blank placeholder paragraphs before the final paragraph.
for ($i = 0; $i -lt $paraXml.Count; $i++) {
  $anchor.InsertParagraphBefore()
}

# Fill each placeholder paragraph with its target content/formatting.
for ($i = 0; $i -lt $paraXml.Count; $i++) {
  $d.Paragraphs($lastIdx + $i).Range.InsertXML($paraXml[$i])
}

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
